# edit.ps1
# Applies the commit "Update gh-pages to output generated at 456a3b4" changes:
#  - Sheet "展览" (Exhibition, sheet 1): bump several "want to go" (F) counters.
#  - Sheet "演出" (Performance, sheet 2): bump two F counters.
#  - Sheet "本地生活" (Local life, sheet 3): the oldest entry
#      (2024-09-06 "HUNTER x HUNTER x animate cafe") expired and was dropped;
#      every later row moved up one slot (keeping its original row-index in
#      column A) and the final (now-duplicated) last row is removed, shrinking
#      the sheet from 15 to 14 data rows. A few F counters also ticked up.
#  - Sheet "全部类型" (All types, sheet 4): mirrors the same F counters.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------------------
# Sheet 1 "展览" - simple counter bumps
# ---------------------------------------------------------------------------
$ws1.Range("F3").Value2  = 2414
$ws1.Range("F6").Value2  = 65
$ws1.Range("F9").Value2  = 2908
$ws1.Range("F11").Value2 = 1046
$ws1.Range("F12").Value2 = 854
$ws1.Range("F15").Value2 = 1483
$ws1.Range("F16").Value2 = 743
$ws1.Range("F17").Value2 = 1703
$ws1.Range("F19").Value2 = 372
$ws1.Range("F21").Value2 = 111
$ws1.Range("F23").Value2 = 2627

# ---------------------------------------------------------------------------
# Sheet 2 "演出" - simple counter bumps
# ---------------------------------------------------------------------------
$ws2.Range("F10").Value2 = 76
$ws2.Range("F38").Value2 = 349

# ---------------------------------------------------------------------------
# Sheet 3 "本地生活" - counter bump on row 6
# ---------------------------------------------------------------------------
$ws3.Range("F6").Value2 = 2499

# Sheet3 (本地生活) rows 8-14: shift content up (drop old row8 'HUNTER x HUNTER' entry), delete old row15
# -- row 8 --
$c = $ws3.Range("B8")
$c.NumberFormat = "@"
$c.Value2 = '2024-09-09'
$ws3.Range("C8").Value2 = '上海·日漫咖啡体验'
$ws3.Range("D8").Value2 = '虹桥路1438号高岛屋百货6楼 Oasis漫画喫茶'
$ws3.Range("E8").Value2 = '2024.09.09 10:00-12.31 22:00'
$ws3.Range("F8").Value2 = 134
$ws3.Range("G8").Value2 = 60
$ws3.Range("H8").Value2 = 'https://show.bilibili.com/platform/detail.html?id=91993'
$ws3.Range("I8").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/IV5rInWT1725347808557.jpeg'

# -- row 9 --
$c = $ws3.Range("B9")
$c.NumberFormat = "@"
$c.Value2 = '2024-09-10'
$ws3.Range("C9").Value2 = '上海·迷你四驱车赛场'
$ws3.Range("D9").Value2 = '虹桥路1438号高岛屋百货6楼 Oasis漫画喫茶'
$ws3.Range("E9").Value2 = '2024.09.10 10:00-12.31 22:00'
$ws3.Range("F9").Value2 = 6
$ws3.Range("G9").Value2 = 48
$ws3.Range("H9").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92042'
$ws3.Range("I9").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/LzFT5TMO1725348229429.png'

# -- row 10 --
$c = $ws3.Range("B10")
$c.NumberFormat = "@"
$c.Value2 = '2024-09-24'
$ws3.Range("C10").Value2 = '上海·星零界·社交游乐·休闲运动·潮玩派对'
$ws3.Range("D10").Value2 = '长宁路1191号长宁来福士B1 长宁来福士'
$ws3.Range("E10").Value2 = '2024.09.24 10:00-12.31 22:00'
$ws3.Range("F10").Value2 = 15
$ws3.Range("G10").Value2 = 68
$ws3.Range("H10").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92659'
$ws3.Range("I10").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/PHS8s1lu1726221065737.png'

# -- row 11 --
$c = $ws3.Range("B11")
$c.NumberFormat = "@"
$c.Value2 = '2024-09-28'
$ws3.Range("C11").Value2 = '上海·［咒术回战 2024 剧场版 咒术回战 0］主题咖啡厅'
$ws3.Range("D11").Value2 = '大悦城 次元波板糖'
$ws3.Range("E11").Value2 = '2024.09.28 00:00-10.27 23:59'
$ws3.Range("F11").Value2 = 359
$ws3.Range("G11").Value2 = 30
$ws3.Range("H11").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92608'
$ws3.Range("I11").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/DBTiL9sY1726727259104.png'

# -- row 12 --
$c = $ws3.Range("B12")
$c.NumberFormat = "@"
$c.Value2 = '2024-10-01'
$ws3.Range("C12").Value2 = '上海·2024·《世界之外》x  萌果酱谷子咖啡'
$ws3.Range("D12").Value2 = '南京东路340号百联ZX 萌果酱谷子咖啡（百联）'
$ws3.Range("E12").Value2 = '2024.10.01 00:00-12.11 23:59'
$ws3.Range("F12").Value2 = 2791
$ws3.Range("G12").Value2 = 30
$ws3.Range("H12").Value2 = 'https://show.bilibili.com/platform/detail.html?id=93006'
$ws3.Range("I12").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/qtffZOKB1727426243733.png'

# -- row 13 --
$c = $ws3.Range("B13")
$c.NumberFormat = "@"
$c.Value2 = '2024-10-01'
$ws3.Range("C13").Value2 = '上海·三丽鸥家族Sanrio Characters主题餐厅·海滩奇遇季'
$ws3.Range("D13").Value2 = '南京东路800号4楼 上海市第一百货商店-C馆'
$ws3.Range("E13").Value2 = '2024.10.01 00:00-10.18 23:59'
$ws3.Range("F13").Value2 = 361
$ws3.Range("G13").Value2 = 10
$ws3.Range("H13").Value2 = 'https://show.bilibili.com/platform/detail.html?id=93078'
$ws3.Range("I13").Value2 = '//i0.hdslb.com/bfs/openplatform/202409/aiu4g5K21727677592777.png'

# -- row 14 --
$c = $ws3.Range("B14")
$c.NumberFormat = "@"
$c.Value2 = '2024-10-10'
$ws3.Range("C14").Value2 = '上海·「火影忍者疾风传 × animate cafe」'
$ws3.Range("D14").Value2 = '西藏北路198号大悦城北座8楼N809-1 animate cafe上海店'
$ws3.Range("E14").Value2 = '2024.10.10 00:00-11.12 23:59'
$ws3.Range("F14").Value2 = 670
$ws3.Range("G14").Value2 = 30
$ws3.Range("H14").Value2 = 'https://show.bilibili.com/platform/detail.html?id=92883'
$ws3.Range("I14").Value2 = '//i2.hdslb.com/bfs/openplatform/202409/aQIhaIgt1727249498713.png'


# Remove the now-duplicated last row (old row 15's content has already been
# copied up into row 14 above); this also shrinks the sheet dimension from
# A1:I15 to A1:I14, matching the target.
$ws3.Rows.Item(15).Delete()

# ---------------------------------------------------------------------------
# Sheet 4 "全部类型" - mirrors the same counters as sheets 1-3
# ---------------------------------------------------------------------------
$ws4.Range("F7").Value2  = 2414
$ws4.Range("F8").Value2  = 2791
$ws4.Range("F11").Value2 = 670
$ws4.Range("F16").Value2 = 65
$ws4.Range("F20").Value2 = 1046
$ws4.Range("F21").Value2 = 854
$ws4.Range("F28").Value2 = 743
$ws4.Range("F31").Value2 = 1703
$ws4.Range("F32").Value2 = 372
$ws4.Range("F42").Value2 = 349
$ws4.Range("F43").Value2 = 2627
